$d = $word.ActiveDocument

# --- Locate the "4. Implement like button" paragraph (the one we are
#     rewriting into "3. stop_following may not work") and the
#     "_GoBack" bookmark that currently sits at the top of the
#     preceding "2. Follows not used in app.py" paragraph. ---

$findRange = $d.Content
$found = $findRange.Find.Execute("4. Implement like button", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetPara = $findRange.Paragraphs(1)
$paraRange = $targetPara.Range
$paraStart = $paraRange.Start
$paraEnd = $paraRange.End

# Replace the whole paragraph's content (minus its trailing paragraph
# mark) with the new wording. A trailing sentinel character is kept
# temporarily so the bookmark can be re-anchored exactly after "work"
# without landing on the paragraph's very-last-character offset.
$bodyRange = $d.Range($paraStart, $paraEnd - 1)
$newText = "3. stop_following may not work"
$bodyRange.Text = $newText + "|"

$anchorPos = $paraStart + $newText.Length

# Move the "_GoBack" bookmark from the previous paragraph to sit right
# after the new paragraph's text (collapsed range).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Drop the sentinel placeholder character now that the bookmark is set.
$sentinelRange = $d.Range($anchorPos, $anchorPos + 1)
$sentinelRange.Text = ""
